# Add 2022-Q3 data: insert a new "2022-Q3" sheet (fund holdings detail) right
# before the existing "2022-Q2" sheet, and update the "总计" (totals) sheet
# with a new top row for the 2022-Q3 totals (pushing the old 2022-Q2 totals
# row down).

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) 总计 sheet: shift the existing 2022-Q2 total row down to row 3, and
#    write the new 2022-Q3 totals into row 2.
# ---------------------------------------------------------------------------

# Carry A2's cell format down to A3 before overwriting A2's neighbours.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.19

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet just before the "2022-Q2" sheet.
#    Duplicating 总计 keeps its header cell style (the one the new sheet's
#    header row ends up using), then we overwrite its contents.
# ---------------------------------------------------------------------------

$totalSheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Header row (B1:D1 already carries the copied header style; extend it to
# E1:H1 so the whole header row matches).
$q3Sheet.Range("B1").Copy()
$q3Sheet.Range("E1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Columns that hold numeric-looking text (fund codes / formatted numbers)
# need to be forced to text so leading zeros / trailing zeros survive.
$q3Sheet.Range("B2:B3").NumberFormat = "@"
$q3Sheet.Range("D2:G3").NumberFormat = "@"

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "013680"
$q3Sheet.Range("C2").Value = "华安品质甄选混合A"
$q3Sheet.Range("D2").Value = "10.86"
$q3Sheet.Range("E2").Value = "42.70"
$q3Sheet.Range("F2").Value = "1.24"
$q3Sheet.Range("G2").Value = "0.1347"
$q3Sheet.Range("H2").Value = 7

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "013681"
$q3Sheet.Range("C3").Value = "华安品质甄选混合C"
$q3Sheet.Range("D3").Value = "4.45"
$q3Sheet.Range("E3").Value = "42.70"
$q3Sheet.Range("F3").Value = "1.24"
$q3Sheet.Range("G3").Value = "0.0552"
$q3Sheet.Range("H3").Value = 7

# Keep the workbook's active tab as it was (总计, the first sheet) rather
# than leaving the newly duplicated sheet active.
$totalSheet.Select()
